# Modal info for rocks/minerals art update
# Insert a new row for "collections"/"COLLECTIONS" key/value pair
# right before the "igneous"/"Igneous" row (currently row 24),
# pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 24 (shifts existing row 24 and below down by one)
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new key/value pair
$ws.Cells.Item(24, 1).Value = "collections"
$ws.Cells.Item(24, 2).Value = "COLLECTIONS"
$ws.Cells.Item(24, 3).Value = 2

# Update the view state to match the saved workbook
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C25").Select()
